$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "shaft" items to more specific real-world-order descriptions.
# Set A10 (Wheel D shaft) before A9 (Central shaft) so the new shared
# strings are appended to the table in that order.
$ws.Cells.Item(10, 1).Value = "Wheel D shaft (5x138mm)"
$ws.Range("A9").Value = "Central shaft (5mm diameter, 100mm long)"

# Row 5: ESC cost updated to real order price.
$ws.Range("C5").Formula = "=16.62*1.101"

# Row 19: PCB cost updated to real order price/formula.
$ws.Range("C19").Formula = "=123.93 + 3.58"

# Row 20: BaneBots Wheels cost becomes a flat literal order price, and the
# quantity becomes a computed fraction instead of a flat literal.
$ws.Range("C20").Value = 74.8
$ws.Range("D20").Formula = "=4/20"

# Move the active selection to E19 to match the saved view state.
[void]$ws.Range("E19").Select()
